{"js": "// Append the continuation sentence about motor creation right after the\n// existing \"Next, to add motors ... yet again and \" sentence, in the same\n// paragraph (matches the author's edit: a second run carrying the new\n// text immediately follows the first, unchanged, run).\n\nconst anchorText =\n  \"Next, to add motors you will select add a device yet again and \";\nconst continuation =\n  \"then select, motor, a port, then a cartridge as previously referenced. \" +\n  \"Repeat this motor creation process unless you have a set of motors \" +\n  \"that you wish to have work as one. In that case select motorgroup.\";\n\nconst body = context.document.body;\nconst results = body.search(anchorText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to extend.\");\n}\n\n// Insert immediately after the matched (unmodified) text, within the same\n// paragraph, so the sentence reads as one continuous flow of text.\nconst target = results.items[0];\ntarget.insertText(continuation, \"After\");\nawait context.sync();\n", "ps1": "# Append the continuation sentence about motor creation right after the\n# existing \"Next, to add motors ... yet again and \" sentence, in the same\n# paragraph (matches the author's edit: a second run carrying the new\n# text immediately follows the first, unchanged, run).\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Next, to add motors you will select add a device yet again and \"\n$continuation = \"then select, motor, a port, then a cartridge as previously referenced. \" + `\n  \"Repeat this motor creation process unless you have a set of motors that you wish to have work as one. \" + `\n  \"In that case select motorgroup.\"\n\n$r = $d.Content\n$r.Find.MatchCase = $true\n$found = $r.Find.Execute($anchorText)\n\nif (-not $found) {\n  throw \"Could not find the target sentence to extend.\"\n}\n\n# Collapse to the end of the matched text, then insert the continuation\n# immediately after it so the sentence reads as one continuous flow.\n$r.Collapse(0)\n$r.InsertAfter($continuation)\n"}
